$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 587, pushing the existing rows 587..653 down to 588..654.
$ws.Rows.Item(587).Insert()

# Fill the new row 587 with the latest weekly record (same constant columns as
# the rest of this "Femacal de La Calera - Cilantro" block; only the
# date/volume/price columns are new).
$ws.Cells.Item(587, 1).Value = 3
$ws.Cells.Item(587, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(587, 3).Value = "Coquimbo"
$ws.Cells.Item(587, 4).Value = 45212
$ws.Cells.Item(587, 4).NumberFormat = $ws.Cells.Item(588, 4).NumberFormat
$ws.Cells.Item(587, 5).Value = 5
$ws.Cells.Item(587, 6).Value = 100112040
$ws.Cells.Item(587, 7).Value = "Cilantro"
$ws.Cells.Item(587, 8).Value = "Sin especificar"
$ws.Cells.Item(587, 9).Value = "Primera"
$ws.Cells.Item(587, 10).Value = 120
$ws.Cells.Item(587, 11).Value = 4000
$ws.Cells.Item(587, 12).Value = 4000
$ws.Cells.Item(587, 13).Value = 4000
$ws.Cells.Item(587, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(587, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(587, 16).Value = 1333
$ws.Cells.Item(587, 17).Value = 3
$ws.Cells.Item(587, 18).Value = "Hortaliza"
